$wb = $excel.ActiveWorkbook
$sheetBefore = $wb.Worksheets.Item("تولید نفت")
$newSheet = $wb.Worksheets.Add($sheetBefore)
$newSheet.Name = "هزینه ها"

$newSheet.Range("A1").Value = "سال"
$newSheet.Range("B1").Value = "امور خدمات عمومی "
$newSheet.Range("C1").Value = "امور دفاعی"
$newSheet.Range("D1").Value = "امور قضایی"
$newSheet.Range("E1").Value = "امور اقتصادی"
$newSheet.Range("F1").Value = "امور محیط زیست"
$newSheet.Range("G1").Value = "امور مسکن"
$newSheet.Range("H1").Value = "امور سلامت"
$newSheet.Range("I1").Value = "امور فرهنگ و تربیت بدنی"
$newSheet.Range("J1").Value = "امور آموزش و پرورش"
$newSheet.Range("K1").Value = "امور رفاه اجتماعی"
$newSheet.Range("L1").Value = "جمع"

$newSheet.Range("A2").Value = 1398
$newSheet.Range("A3").Value = 1399
$newSheet.Range("A4").Value = 1400
$newSheet.Range("A5").Value = 1401
$newSheet.Range("A6").Value = 1402

$newSheet.Range("B5").Value = 318045716
$newSheet.Range("C5").Value = 1951399111
$newSheet.Range("D5").Value = 292308709
$newSheet.Range("E5").Value = 198775474
$newSheet.Range("F5").Value = 15164449
$newSheet.Range("G5").Value = 3190717
$newSheet.Range("H5").Value = 1207892378
$newSheet.Range("I5").Value = 156162709
$newSheet.Range("J5").Value = 1790139379
$newSheet.Range("K5").Value = 2833391544
$newSheet.Range("L5").Formula = "=SUM(B5:K5)"

$newSheet.Range("B6").Value = 381179603
$newSheet.Range("C6").Value = 2668999100
$newSheet.Range("D6").Value = 415705303
$newSheet.Range("E6").Value = 244325051
$newSheet.Range("F6").Value = 25503187
$newSheet.Range("G6").Value = 5651500
$newSheet.Range("H6").Value = 1816484885
$newSheet.Range("I6").Value = 215842143
$newSheet.Range("J6").Value = 2843718967
$newSheet.Range("K6").Value = 5429761932
$newSheet.Range("L6").Formula = "=SUM(B6:K6)"

Write-Output "done"
